$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the test-method data in rows 2 and 3 (columns A, E, F) ---
$scratch = $ws.Range("H1")
foreach ($col in @("A","E","F")) {
    $c2 = $ws.Range("$col`2")
    $c3 = $ws.Range("$col`3")
    $c2.Copy($scratch)
    $c3.Copy($c2)
    $scratch.Copy($c3)
}
$scratch.Clear()

# --- Swap which hyperlink relationship (rId) is used by B2 / B3 ---
# Preserve original formatting of the hyperlink cells first.
$ws.Range("B2").Copy($ws.Range("H1"))

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "https://uat.royalbank.can-act.com/")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://uat.royalbank.can-act.com/")

$ws.Range("H1").Copy($ws.Range("B2"))
$ws.Range("H1").Copy($ws.Range("B3"))
$ws.Range("H1").Clear()

# --- Reset the view: scroll back to A1, drop the stale selection on E3 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
